$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet "question_answers": update answer values in column B (text cells) ---
$ws1.Cells.Item(2,2).NumberFormat = "@"
$ws1.Cells.Item(2,2).Value = "1"
$ws1.Cells.Item(3,2).NumberFormat = "@"
$ws1.Cells.Item(3,2).Value = "4"
$ws1.Cells.Item(5,2).NumberFormat = "@"
$ws1.Cells.Item(5,2).Value = "3"
$ws1.Cells.Item(6,2).NumberFormat = "@"
$ws1.Cells.Item(6,2).Value = "4"
$ws1.Cells.Item(7,2).NumberFormat = "@"
$ws1.Cells.Item(7,2).Value = "2"
$ws1.Cells.Item(9,2).NumberFormat = "@"
$ws1.Cells.Item(9,2).Value = "3"
$ws1.Cells.Item(10,2).NumberFormat = "@"
$ws1.Cells.Item(10,2).Value = "1"
$ws1.Cells.Item(11,2).NumberFormat = "@"
$ws1.Cells.Item(11,2).Value = "4"
$ws1.Cells.Item(12,2).NumberFormat = "@"
$ws1.Cells.Item(12,2).Value = "1"
$ws1.Cells.Item(13,2).NumberFormat = "@"
$ws1.Cells.Item(13,2).Value = "2"
$ws1.Cells.Item(14,2).NumberFormat = "@"
$ws1.Cells.Item(14,2).Value = "3"
$ws1.Cells.Item(15,2).NumberFormat = "@"
$ws1.Cells.Item(15,2).Value = "1"
$ws1.Cells.Item(16,2).NumberFormat = "@"
$ws1.Cells.Item(16,2).Value = "3"
$ws1.Cells.Item(17,2).NumberFormat = "@"
$ws1.Cells.Item(17,2).Value = "1"
$ws1.Cells.Item(18,2).NumberFormat = "@"
$ws1.Cells.Item(18,2).Value = "1"
$ws1.Cells.Item(19,2).NumberFormat = "@"
$ws1.Cells.Item(19,2).Value = "3"
$ws1.Cells.Item(20,2).NumberFormat = "@"
$ws1.Cells.Item(20,2).Value = "2"
$ws1.Cells.Item(21,2).NumberFormat = "@"
$ws1.Cells.Item(21,2).Value = "3"
$ws1.Cells.Item(22,2).NumberFormat = "@"
$ws1.Cells.Item(22,2).Value = "2"
$ws1.Cells.Item(23,2).NumberFormat = "@"
$ws1.Cells.Item(23,2).Value = "4"
$ws1.Cells.Item(25,2).NumberFormat = "@"
$ws1.Cells.Item(25,2).Value = "4"
$ws1.Cells.Item(26,2).NumberFormat = "@"
$ws1.Cells.Item(26,2).Value = "4"
$ws1.Cells.Item(28,2).NumberFormat = "@"
$ws1.Cells.Item(28,2).Value = "4"
$ws1.Cells.Item(29,2).NumberFormat = "@"
$ws1.Cells.Item(29,2).Value = "2"
$ws1.Cells.Item(31,2).NumberFormat = "@"
$ws1.Cells.Item(31,2).Value = "3"
$ws1.Cells.Item(32,2).NumberFormat = "@"
$ws1.Cells.Item(32,2).Value = "2"
$ws1.Cells.Item(35,2).NumberFormat = "@"
$ws1.Cells.Item(35,2).Value = "1"
$ws1.Cells.Item(36,2).NumberFormat = "@"
$ws1.Cells.Item(36,2).Value = "2"
$ws1.Cells.Item(37,2).NumberFormat = "@"
$ws1.Cells.Item(37,2).Value = "3"
$ws1.Cells.Item(39,2).NumberFormat = "@"
$ws1.Cells.Item(39,2).Value = "4"
$ws1.Cells.Item(41,2).NumberFormat = "@"
$ws1.Cells.Item(41,2).Value = "1"
$ws1.Cells.Item(42,2).NumberFormat = "@"
$ws1.Cells.Item(42,2).Value = "3"
$ws1.Cells.Item(44,2).NumberFormat = "@"
$ws1.Cells.Item(44,2).Value = "4"
$ws1.Cells.Item(48,2).NumberFormat = "@"
$ws1.Cells.Item(48,2).Value = "1"
$ws1.Cells.Item(49,2).NumberFormat = "@"
$ws1.Cells.Item(49,2).Value = "1"
$ws1.Cells.Item(50,2).NumberFormat = "@"
$ws1.Cells.Item(50,2).Value = "4"
$ws1.Cells.Item(51,2).NumberFormat = "@"
$ws1.Cells.Item(51,2).Value = "1"
$ws1.Cells.Item(52,2).NumberFormat = "@"
$ws1.Cells.Item(52,2).Value = "2"
$ws1.Cells.Item(53,2).NumberFormat = "@"
$ws1.Cells.Item(53,2).Value = "3"
$ws1.Cells.Item(54,2).NumberFormat = "@"
$ws1.Cells.Item(54,2).Value = "1"
$ws1.Cells.Item(55,2).NumberFormat = "@"
$ws1.Cells.Item(55,2).Value = "2"
$ws1.Cells.Item(56,2).NumberFormat = "@"
$ws1.Cells.Item(56,2).Value = "1"
$ws1.Cells.Item(57,2).NumberFormat = "@"
$ws1.Cells.Item(57,2).Value = "1"
$ws1.Cells.Item(58,2).NumberFormat = "@"
$ws1.Cells.Item(58,2).Value = "2"
$ws1.Cells.Item(59,2).NumberFormat = "@"
$ws1.Cells.Item(59,2).Value = "1"
$ws1.Cells.Item(60,2).NumberFormat = "@"
$ws1.Cells.Item(60,2).Value = "3"
$ws1.Cells.Item(61,2).NumberFormat = "@"
$ws1.Cells.Item(61,2).Value = "4"
$ws1.Cells.Item(62,2).NumberFormat = "@"
$ws1.Cells.Item(62,2).Value = "1"
$ws1.Cells.Item(63,2).NumberFormat = "@"
$ws1.Cells.Item(63,2).Value = "4"
$ws1.Cells.Item(65,2).NumberFormat = "@"
$ws1.Cells.Item(65,2).Value = "1"
$ws1.Cells.Item(69,2).NumberFormat = "@"
$ws1.Cells.Item(69,2).Value = "2"
$ws1.Cells.Item(73,2).NumberFormat = "@"
$ws1.Cells.Item(73,2).Value = "4"
$ws1.Cells.Item(75,2).NumberFormat = "@"
$ws1.Cells.Item(75,2).Value = "3"
$ws1.Cells.Item(76,2).NumberFormat = "@"
$ws1.Cells.Item(76,2).Value = "2"
$ws1.Cells.Item(77,2).NumberFormat = "@"
$ws1.Cells.Item(77,2).Value = "2"
$ws1.Cells.Item(78,2).NumberFormat = "@"
$ws1.Cells.Item(78,2).Value = "2"
$ws1.Cells.Item(79,2).NumberFormat = "@"
$ws1.Cells.Item(79,2).Value = "2"
$ws1.Cells.Item(80,2).NumberFormat = "@"
$ws1.Cells.Item(80,2).Value = "1"
$ws1.Cells.Item(81,2).NumberFormat = "@"
$ws1.Cells.Item(81,2).Value = "2"
$ws1.Cells.Item(82,2).NumberFormat = "@"
$ws1.Cells.Item(82,2).Value = "3"
$ws1.Cells.Item(83,2).NumberFormat = "@"
$ws1.Cells.Item(83,2).Value = "4"
$ws1.Cells.Item(85,2).NumberFormat = "@"
$ws1.Cells.Item(85,2).Value = "1"
$ws1.Cells.Item(86,2).NumberFormat = "@"
$ws1.Cells.Item(86,2).Value = "1"
$ws1.Cells.Item(87,2).NumberFormat = "@"
$ws1.Cells.Item(87,2).Value = "2"
$ws1.Cells.Item(88,2).NumberFormat = "@"
$ws1.Cells.Item(88,2).Value = "4"
$ws1.Cells.Item(89,2).NumberFormat = "@"
$ws1.Cells.Item(89,2).Value = "2"
$ws1.Cells.Item(90,2).NumberFormat = "@"
$ws1.Cells.Item(90,2).Value = "1"
$ws1.Cells.Item(91,2).NumberFormat = "@"
$ws1.Cells.Item(91,2).Value = "3"
$ws1.Cells.Item(92,2).NumberFormat = "@"
$ws1.Cells.Item(92,2).Value = "4"
$ws1.Cells.Item(94,2).NumberFormat = "@"
$ws1.Cells.Item(94,2).Value = "4"
$ws1.Cells.Item(95,2).NumberFormat = "@"
$ws1.Cells.Item(95,2).Value = "1"
$ws1.Cells.Item(96,2).NumberFormat = "@"
$ws1.Cells.Item(96,2).Value = "4"
$ws1.Cells.Item(97,2).NumberFormat = "@"
$ws1.Cells.Item(97,2).Value = "3"
$ws1.Cells.Item(98,2).NumberFormat = "@"
$ws1.Cells.Item(98,2).Value = "4"
$ws1.Cells.Item(99,2).NumberFormat = "@"
$ws1.Cells.Item(99,2).Value = "2"
$ws1.Cells.Item(100,2).NumberFormat = "@"
$ws1.Cells.Item(100,2).Value = "1"
$ws1.Cells.Item(102,2).NumberFormat = "@"
$ws1.Cells.Item(102,2).Value = "2"
$ws1.Cells.Item(103,2).NumberFormat = "@"
$ws1.Cells.Item(103,2).Value = "3"
$ws1.Cells.Item(105,2).NumberFormat = "@"
$ws1.Cells.Item(105,2).Value = "2"
$ws1.Cells.Item(106,2).NumberFormat = "@"
$ws1.Cells.Item(106,2).Value = "4"
$ws1.Cells.Item(107,2).NumberFormat = "@"
$ws1.Cells.Item(107,2).Value = "3"
$ws1.Cells.Item(108,2).NumberFormat = "@"
$ws1.Cells.Item(108,2).Value = "2"
$ws1.Cells.Item(109,2).NumberFormat = "@"
$ws1.Cells.Item(109,2).Value = "3"
$ws1.Cells.Item(110,2).NumberFormat = "@"
$ws1.Cells.Item(110,2).Value = "2"
$ws1.Cells.Item(111,2).NumberFormat = "@"
$ws1.Cells.Item(111,2).Value = "4"

# --- Sheet "outputs": fix row labels in column A (shifted categories) ---
$ws2.Cells.Item(21,1).Value = "dysthymia_type_a"
$ws2.Cells.Item(22,1).Value = "dysthymia_type_b"
$ws2.Cells.Item(23,1).Value = "dysthymia_total"
$ws2.Cells.Item(24,1).Value = "autistic_disorder_type_a"
$ws2.Cells.Item(25,1).Value = "autistic_disorder_type_b"
$ws2.Cells.Item(26,1).Value = "autistic_disorder_type_c"
$ws2.Cells.Item(27,1).Value = "autistic_disorder_total"
$ws2.Cells.Item(28,1).Value = "asperger_disorder_type_a"
$ws2.Cells.Item(29,1).Value = "asperger_disorder_type_b"
$ws2.Cells.Item(30,1).Value = "asperger_disorder_total"
$ws2.Cells.Item(31,1).Value = "social_phobia"
$ws2.Cells.Item(32,1).Value = "seperation_anxiety_disorder"
$ws2.Cells.Item(33,1).Value = "enuresis"

# --- Sheet "outputs": update computed score values in column B (numeric cells) ---
$ws2.Cells.Item(2,2).Value = 7
$ws2.Cells.Item(3,2).Value = 5
$ws2.Cells.Item(4,2).Value = 12
$ws2.Cells.Item(5,2).Value = 8
$ws2.Cells.Item(6,2).Value = 15
$ws2.Cells.Item(8,2).Value = 1
$ws2.Cells.Item(9,2).Value = 4
$ws2.Cells.Item(10,2).Value = 5
$ws2.Cells.Item(17,2).Value = 2
$ws2.Cells.Item(18,2).Value = 2
$ws2.Cells.Item(20,2).Value = 4
$ws2.Cells.Item(21,2).Value = 1
$ws2.Cells.Item(22,2).Value = 2
$ws2.Cells.Item(23,2).Value = 3
$ws2.Cells.Item(24,2).Value = 3
$ws2.Cells.Item(25,2).Value = 4
$ws2.Cells.Item(27,2).Value = 10
$ws2.Cells.Item(28,2).Value = 3
$ws2.Cells.Item(29,2).Value = 3
$ws2.Cells.Item(31,2).Value = 2
$ws2.Cells.Item(32,2).Value = 7
$ws2.Cells.Item(33,2).Value = 1
$ws2.Cells.Item(37,2).Value = 1
$ws2.Cells.Item(38,2).Value = 0
$ws2.Cells.Item(44,2).Value = 1
$ws2.Cells.Item(46,2).Value = 1
$ws2.Cells.Item(48,2).Value = 1
